$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.394.94"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3
$ws.Range("D3").Value = "1.848.51"
$ws.Range("E3").Value = "  +0.07%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.51"
$ws.Range("E5").Value = "  +0.10%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6296"
$ws.Range("E6").Value = "  -0.43%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07622"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2940"
$ws.Range("E9").Value = "  -0.73%  "

# Row 10
$ws.Range("E10").Value = "  +0.25%  "

# Row 12
$ws.Range("D12").Value = "1.834.02"
$ws.Range("E12").Value = "  -0.58%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.008"
$ws.Range("E13").Value = "  +0.23%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001089"
$ws.Range("E14").Value = "  +7.75%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6794"
$ws.Range("E15").Value = "  -0.88%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.59"

# Row 17
$ws.Range("D17").Value = "2.078.46"
$ws.Range("E17").Value = "  -8.21%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.152"
$ws.Range("E18").Value = "  -0.09%  "

# Row 19
$ws.Range("D19").Value = "29.421.77"
$ws.Range("E19").Value = "  -0.01%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.73"
$ws.Range("E20").Value = "  -0.49%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.44"
$ws.Range("E21").Value = "  -0.05%  "

# Row 22
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.461"
$ws.Range("E23").Value = "  -1.37%  "

# Row 24
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.79"
$ws.Range("E25").Value = "  +0.60%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1390"
$ws.Range("E26").Value = "  -0.80%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.381"
$ws.Range("E27").Value = "  -0.08%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.64"
$ws.Range("E28").Value = "  -0.08%  "

# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.466"
$ws.Range("E29").Value = "  +0.17%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.309"
$ws.Range("E30").Value = "  +4.65%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05647"
$ws.Range("E31").Value = "  -1.52%  "

# Row 32
$ws.Range("E32").Value = "  -0.48%  "

# Row 33
$ws.Range("E33").Value = "  +0.23%  "

# Row 34
$ws.Range("E34").Value = "  +0.04%  "

# Row 35
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7095"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.586"
$ws.Range("E37").Value = "  -0.15%  "

# Row 38
$ws.Range("D38").Value = "1.233.92"
$ws.Range("E38").Value = "  -1.59%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.773"
$ws.Range("E39").Value = "  -0.23%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01797"
$ws.Range("E40").Value = "  -1.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.459"
$ws.Range("E41").Value = "  +4.57%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9095"
$ws.Range("E42").Value = "  -0.10%  "

# Row 44
$ws.Range("D44").Value = "1.988.37"
$ws.Range("E44").Value = "  -0.68%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.46"
$ws.Range("E45").Value = "  -0.24%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.00"
$ws.Range("E46").Value = "  -0.34%  "

# Row 47
$ws.Range("E47").Value = "  +3.13%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.161"
$ws.Range("E48").Value = "  +1.34%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4017"
$ws.Range("E49").Value = "  -0.28%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.013"
$ws.Range("E50").Value = "  -1.64%  "

# Row 51
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.686"
$ws.Range("E51").Value = "  -1.04%  "
